$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "E"=3; "F"=1; "G"=2.651738666666667; "H"=7.955216; "I"=0.5604432756981275; "J"=0.5604432756981274; "K"=3; "L"=1; "M"=11.23020933333333; "N"=33.690628; "O"=0.4111393124289954; "P"=0.4111393124289954; "Q"=29.77958032396089; "R"=268.016222915648; "S"=0.230420263025982; "T"=0.230420263025982 }
  3 = @{ "E"=3; "F"=1; "G"=2.651738666666667; "H"=7.955216; "I"=0.5604432756981275; "J"=0.5604432756981274; "K"=3; "L"=1; "M"=16.08464233333333; "N"=48.253927; "O"=0.5888606875710045; "P"=0.5888606875710045; "Q"=42.65226801480356; "R"=383.870412133232; "S"=0.3300230126721454; "T"=0.3300230126721453 }
  4 = @{ "E"=3; "F"=1; "G"=0.5155733333333333; "H"=1.54672; "I"=0.10896609512398; "J"=0.10896609512398; "K"=3; "L"=1; "M"=11.23020933333333; "N"=33.690628; "O"=0.4111393124289954; "P"=0.4111393124289954; "Q"=5.789996460017777; "R"=52.10996814015999; "S"=0.04480024542734565; "T"=0.04480024542734563 }
  5 = @{ "E"=3; "F"=1; "G"=0.5155733333333333; "H"=1.54672; "I"=0.10896609512398; "J"=0.10896609512398; "K"=3; "L"=1; "M"=16.08464233333333; "N"=48.253927; "O"=0.5888606875710045; "P"=0.5888606875710045; "Q"=8.292812663271112; "R"=74.63531396944; "S"=0.06416584969663434; "T"=0.06416584969663433 }
  6 = @{ "E"=3; "F"=1; "G"=1.564190333333333; "H"=4.692571; "I"=0.3305906291778926; "J"=0.3305906291778925; "K"=3; "L"=1; "M"=11.23020933333333; "N"=33.690628; "O"=0.4111393124289954; "P"=0.4111393124289954; "Q"=17.56618488050978; "R"=158.095663924588; "S"=0.1359188039756677; "T"=0.1359188039756677 }
  7 = @{ "E"=3; "F"=1; "G"=1.564190333333333; "H"=4.692571; "I"=0.3305906291778926; "J"=0.3305906291778925; "K"=3; "L"=1; "M"=16.08464233333333; "N"=48.253927; "O"=0.5888606875710045; "P"=0.5888606875710045; "Q"=25.15944205292411; "R"=226.434978476317; "S"=0.1946718252022248; "T"=0.1946718252022248 }
}

foreach ($row in $data.Keys) {
  foreach ($col in $data[$row].Keys) {
    $ws.Range("$col$row").Value = $data[$row][$col]
  }
}
